$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) Split "Jatin Singh Bisht (21105039)" into multiple runs with
#    identical formatting (mirrors the real edit history where the
#    roll number got retyped a few characters at a time).
# -----------------------------------------------------------------
$para = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Jatin Singh Bisht (21105039)*") {
        $para = $p
        break
    }
}

if ($para -ne $null) {
    $pr = $para.Range
    $full = $pr.Text
    $idx = $full.IndexOf(" (21105039)")
    if ($idx -ge 0) {
        $subStart = $pr.Start + $idx
        # stop before the trailing paragraph mark
        $subEnd = $pr.Start + $full.Length - 1
        $subRange = $d.Range($subStart, $subEnd)

        $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
          '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
          '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
          '<pkg:xmlData>' +
          '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body><w:p>' +
          '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> (211050</w:t></w:r>' +
          '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>3</w:t></w:r>' +
          '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>9</w:t></w:r>' +
          '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>)</w:t></w:r>' +
          '</w:p></w:body></w:document>' +
          '</pkg:xmlData></pkg:part></pkg:package>'

        $subRange.InsertXML($xml)
    }
}

# -----------------------------------------------------------------
# 2) Stamp each inline drawing with a wp14:anchorId / wp14:editId
#    pair (Word mints these whenever it resaves a document that
#    contains drawings).
# -----------------------------------------------------------------
$ids = @(
    @{ anchor = "4A56C881"; edit = "7E56CBC2" },
    @{ anchor = "62D3DBF1"; edit = "74A9A02A" },
    @{ anchor = "2A724F83"; edit = "7E612FED" },
    @{ anchor = "7178810B"; edit = "36FB7D34" },
    @{ anchor = "68FF179E"; edit = "1D39FF89" },
    @{ anchor = "550CFD44"; edit = "6E5CF304" },
    @{ anchor = "7EE3CE12"; edit = "04657EDB" },
    @{ anchor = "44C043D7"; edit = "1EF7765C" },
    @{ anchor = "16846670"; edit = "705AA74F" },
    @{ anchor = "60D16C9F"; edit = "5D40F164" }
)

# Range.WordOpenXML re-exports a *self-contained* mini-package, so
# image relationship ids inside it are re-minted starting at rId4 for
# every single shape (they do NOT reflect the real, document-wide
# r:embed ids once there is more than one picture). Pull the whole
# document's OOXML once instead - relationship ids in that export do
# line up with the real word/_rels/document.xml.rels - and slice the
# individual drawing runs out of it.
$fullXml = $d.Content.WordOpenXML
$pattern = '(?s)<w:r>(?:(?!</?w:r>).)*<w:drawing>.*?</w:drawing>(?:(?!</?w:r>).)*</w:r>'

$drawingRuns = @()
$remaining = $fullXml
while ($remaining -match $pattern) {
    $val = $matches[0]
    $drawingRuns += $val
    $idx = $remaining.IndexOf($val)
    $remaining = $remaining.Substring($idx + $val.Length)
}

$count = $d.InlineShapes.Count
for ($i = 1; $i -le $count; $i++) {
    if ($i -gt $drawingRuns.Count) { continue }

    $runXml = $drawingRuns[$i - 1]

    $pair = $ids[$i - 1]
    $anchorId = $pair.anchor
    $editId = $pair.edit

    $newRunXml = $runXml -replace '<wp:inline distT="0" distB="0" distL="0" distR="0">', ('<wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="' + $anchorId + '" wp14:editId="' + $editId + '">')

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
      '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
      '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
      '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture">' +
      '<w:body><w:p>' + $newRunXml + '</w:p></w:body></w:document>' +
      '</pkg:xmlData></pkg:part></pkg:package>'

    # A plain InsertXML on the shape's own Range *inserts* a sibling
    # copy rather than replacing it (the Range collapses to the
    # drawing's anchor character, which cannot be overwritten as text).
    # Delete the shape first, then insert the (attribute-augmented)
    # replacement drawing at the now-collapsed insertion point so the
    # picture count / order / relationship ids stay exactly as before.
    $shp = $d.InlineShapes.Item($i)
    $insStart = $shp.Range.Start
    $null = $shp.Delete()
    $target = $d.Range($insStart, $insStart)
    $null = $target.InsertXML($xml)
}

Write-Host "Done"
